# Updates the cryptocurrency price table (cols D and E) plus a row swap
# (Filecoin <-> WrappedliquidstakedEther2.0) for rows 31/32, matching the
# latest scrape as described in the commit message.

function Set-CellText {
    param($ws, $ref, $val)
    $cell = $ws.Range($ref)
    # Preserve the cell's existing style, but force a text number-format while
    # assigning the value so Excel doesn't auto-convert numeric-looking
    # strings (e.g. "1.001", "302.48") into floating point numbers.
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" "23.521.36"
Set-CellText $ws "E2" "  +1.28%  "
Set-CellText $ws "D3" "1.654.96"
Set-CellText $ws "D4" "1.001"
Set-CellText $ws "E4" "  -0.09%  "
Set-CellText $ws "E5" "  -0.11%  "
Set-CellText $ws "D6" "302.34"
Set-CellText $ws "E6" "  -0.30%  "
Set-CellText $ws "D7" "0.3835"
Set-CellText $ws "E7" "  +1.31%  "
Set-CellText $ws "D8" "51.16"
Set-CellText $ws "E8" "  -0.80%  "
Set-CellText $ws "D9" "0.3592"
Set-CellText $ws "E9" "  +1.71%  "
Set-CellText $ws "D10" "1.242"
Set-CellText $ws "E10" "  +3.08%  "
Set-CellText $ws "D11" "0.08188"
Set-CellText $ws "E11" "  +1.04%  "
Set-CellText $ws "D12" "1.001"
Set-CellText $ws "E12" "  -0.07%  "
Set-CellText $ws "D13" "22.42"
Set-CellText $ws "E13" "  +1.04%  "
Set-CellText $ws "D14" "6.493"
Set-CellText $ws "E14" "  +2.17%  "
Set-CellText $ws "D15" "7.493"
Set-CellText $ws "E16" "  +0.76%  "
Set-CellText $ws "D17" "1.657.60"
Set-CellText $ws "E17" "  +3.30%  "
Set-CellText $ws "D18" "97.48"
Set-CellText $ws "D19" "0.06976"
Set-CellText $ws "E19" "  +1.02%  "
Set-CellText $ws "D20" "6.822"
Set-CellText $ws "E20" "  +5.33%  "
Set-CellText $ws "E21" "  +2.59%  "
Set-CellText $ws "D22" "1.002"
Set-CellText $ws "E22" "  +0.01%  "
Set-CellText $ws "D23" "12.68"
Set-CellText $ws "E23" "  +2.90%  "
Set-CellText $ws "D24" "23.530.59"
Set-CellText $ws "E24" "  +1.37%  "
Set-CellText $ws "D25" "2.502"
Set-CellText $ws "E25" "  -0.44%  "
Set-CellText $ws "D26" "3.001"
Set-CellText $ws "E26" "  -0.73%  "
Set-CellText $ws "D27" "21.22"
Set-CellText $ws "E27" "  +1.69%  "
Set-CellText $ws "D28" "152.00"
Set-CellText $ws "E28" "  +0.53%  "
Set-CellText $ws "D29" "5.242"
Set-CellText $ws "E29" "  +0.12%  "
Set-CellText $ws "D30" "133.86"
Set-CellText $ws "E30" "  +1.36%  "
Set-CellText $ws "B31" "WrappedliquidstakedEther2.0"
Set-CellText $ws "C31" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-CellText $ws "D31" "1.837.72"
Set-CellText $ws "E31" "  +2.90%  "
Set-CellText $ws "B32" "Filecoin"
Set-CellText $ws "C32" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-CellText $ws "D32" "7.198"
Set-CellText $ws "E32" "  +11.43%  "
Set-CellText $ws "D33" "2.250"
Set-CellText $ws "E33" "  +7.37%  "
Set-CellText $ws "D34" "12.09"
Set-CellText $ws "E34" "  +7.04%  "
Set-CellText $ws "E35" "  -0.52%  "
Set-CellText $ws "E36" "  +3.64%  "
Set-CellText $ws "D37" "6.131"
Set-CellText $ws "E37" "  +4.79%  "
Set-CellText $ws "D38" "0.2499"
Set-CellText $ws "E38" "  +1.91%  "
Set-CellText $ws "D39" "0.08792"
Set-CellText $ws "D40" "0.07023"
Set-CellText $ws "E40" "  +1.24%  "
Set-CellText $ws "D41" "13.22"
Set-CellText $ws "E41" "  +10.73%  "
Set-CellText $ws "D42" "0.7015"
Set-CellText $ws "E42" "  +1.96%  "
Set-CellText $ws "E43" "  +0.69%  "
Set-CellText $ws "D44" "15.99"
Set-CellText $ws "E44" "  +4.31%  "
Set-CellText $ws "D45" "0.6531"
Set-CellText $ws "E45" "  +3.34%  "
Set-CellText $ws "D46" "1.000"
Set-CellText $ws "E46" "  -0.03%  "
Set-CellText $ws "D47" "2.308"
Set-CellText $ws "E47" "  +2.50%  "
Set-CellText $ws "D48" "3.957"
Set-CellText $ws "E48" "  +0.29%  "
Set-CellText $ws "D49" "0.07912"
Set-CellText $ws "E49" "  +0.66%  "
Set-CellText $ws "D50" "128.04"
Set-CellText $ws "E50" "  +0.73%  "
Set-CellText $ws "D51" "1.193"
Set-CellText $ws "E51" "  +1.82%  "
